$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in the header row.
$ws.Range("F1").Value = "Last status check on: 22.02.2022 06:00"

# Row 4 (Globus) got a fresh price scrape:
#  - new price goes into B4
#  - the previous B4 price (36.1) becomes the "Old Cena" in C4
#  - the delta is now rendered as a signed text string in D4
#  - the "Old Datum" column now stores a plain text timestamp in E4
$ws.Range("C4").Value = 36.1
$ws.Range("B4").Value = 36.9

# D4 / E4 switch from numeric/date cells to plain text cells. Force text
# storage via NumberFormat "@" so Excel doesn't re-parse the literal back
# into a number/date, then restore the default "Normal" style so no
# lingering number format is left applied to the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "+0.8"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "2022-02-22 06:00:09"
$ws.Range("E4").Style = "Normal"
